$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the match data (columns F:V) between row 60 and row 61 ---
# Row 60 currently holds the "Alashkert vs BKMA" match, row 61 holds the
# "Pyunik Yerevan vs Urartu" match. The author swapped their order while
# keeping the row index (column A) and the shared B:E columns untouched.

$row60 = $ws.Range("F60:V60").Value2
$row61 = $ws.Range("F61:V61").Value2

$ws.Range("F61:V61").Value = $row60
$ws.Range("F60:V60").Value = $row61

# --- Append the new match row (row 79) ---

# Copy number formatting (borders/bold/alignment for A, date format for E)
# from the last existing data row so the new row matches the sheet style.
$ws.Range("A78").Copy()
$ws.Range("A79").PasteSpecial(-4122)
$ws.Range("E78").Copy()
$ws.Range("E79").PasteSpecial(-4122)

$ws.Range("A79").Value = 78
$ws.Range("B79").Value = "armenia"
$ws.Range("C79").Value = "premier-league"
$ws.Range("D79").Value = "2023-2024"
$ws.Range("E79").Value = 45242.47916666666
$ws.Range("F79").Value = "Van"
$ws.Range("G79").Value = 2
$ws.Range("H79").Value = "Urartu"
$ws.Range("I79").Value = 2
$ws.Range("J79").Value = 5.98
$ws.Range("K79").Value = "10/11/2023 23:42"
$ws.Range("L79").Value = 8.449999999999999
$ws.Range("M79").Value = "12/11/2023 11:18"
$ws.Range("N79").Value = 4.64
$ws.Range("O79").Value = "10/11/2023 23:42"
$ws.Range("P79").Value = 5.08
$ws.Range("Q79").Value = "12/11/2023 11:18"
$ws.Range("R79").Value = 1.39
$ws.Range("S79").Value = "10/11/2023 23:42"
$ws.Range("T79").Value = 1.35
$ws.Range("U79").Value = "12/11/2023 11:18"
$ws.Range("V79").Value = "https://www.betexplorer.com/football/armenia/premier-league/van-urartu/02k5wi4B/"
